# Update "Horarios" workbook with the latest scrape for Línea 141 - 190
$wb = $excel.ActiveWorkbook

$newTime = "02:40:13"

# --- Sheet "LP1912" ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"

$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "02:58"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 18

$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "03:56"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 76

$ws1.Range("A8").Value = $newTime
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 81

# --- Sheet "LP1912-215" ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "02:58"
$ws2.Range("D6").Value = 18

# --- Sheet "6203-6173" ---
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
